$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns ---------------------------------------------------
# Original layout: A = labels, B = lowercase "no conformidad" reasons (rows 6-12).
# New layout: A unchanged, a narrow "B" column holding a literal "('",
# the (now capitalized) reason text moves to "C", and a narrow "D" column
# holds the closing "', 1) ," (or "', 1)" on the last row) - these look like
# fragments of generated Python tuple/list source.
#
# Insert a blank column at B; this pushes the existing reason text from B
# into C automatically (and shifts the dimension/spans along with it).
$ws.Columns.Item(2).Insert()

# --- Capitalize the reason text now sitting in column C (rows 6-12) -------
$ws.Range("C6").Value = "Falla en el cumplimiento de las especificaciones de una actividad de capacitacion"
$ws.Range("C7").Value = "Falla en el cumplimiento de disponer la infraestructura y equipos necesarios"
$ws.Range("C8").Value = "Materiales de servicios de capacitacion inadecuados"
$ws.Range("C9").Value = "Falla en el cumplimiento de las exigencias de seguridad y salud ocupacional"
$ws.Range("C10").Value = "Deficiencias en el sistema de gestion de la calidad del organismo de capacitacion"
$ws.Range("C11").Value = "Instrumentos de capacitacion inadecuados"
$ws.Range("C12").Value = "Relatores y/o facilitadores no son evaluados  en terminos de desempeño, dentro de las espectativas del organsmo"

# --- Add the "('" prefix marker cells in column B --------------------------
$ws.Range("B6").Value = "'('"
$ws.Range("B7").Value = "'('"
$ws.Range("B8").Value = "'('"
$ws.Range("B9").Value = "'('"
$ws.Range("B10").Value = "'('"
$ws.Range("B11").Value = "'('"
$ws.Range("B12").Value = "'('"

# --- Add the closing "', 1) ," / "', 1)" marker cells in column D ----------
# NOTE: a leading single-quote in a COM Range.Value is consumed by Excel as
# the "quote prefix" marker (forces text, doesn't become part of the
# stored string) - double it so one literal apostrophe survives as data.
$ws.Range("D6").Value = "'', 1) ,"
$ws.Range("D7").Value = "'', 1) ,"
$ws.Range("D8").Value = "'', 1) ,"
$ws.Range("D9").Value = "'', 1) ,"
$ws.Range("D10").Value = "'', 1) ,"
$ws.Range("D11").Value = "'', 1) ,"
$ws.Range("D12").Value = "'', 1)"

# --- Column widths (best-effort match to the new narrow helper columns) ---
$ws.Range("B1").ColumnWidth = 1.3333333333333333
$ws.Range("D1").ColumnWidth = 4.333333333333333
$ws.Range("E1").ColumnWidth = 28.833333333333332

# --- Selection moves to C3 --------------------------------------------------
$ws.Range("C3").Select()
